{"js": "// Apply strikethrough formatting to three specific \"deliverable\" bullet\n// paragraphs that were marked as done/resolved (data retention policy,\n// write/test triggers, and the \"at least two triggers must be for delete\"\n// bullet). The strike needs to land on every run AND on the paragraph\n// mark itself (so a fresh run typed at the end of the paragraph would\n// inherit it), matching how the other already-struck bullets are built.\n\nconst targets = [\n  \"Data retention policy: data is kept indefinitely. There are numerous ways to resolve this issue, find one and implement it.\",\n  \"Write and test seven (7) triggers for seven (7) separate tables to implement the business rules.\",\n  \"At least two (2) of the triggers must be for delete.\"\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nfor (const paragraph of paragraphs.items) {\n  const text = paragraph.text.trim();\n  if (targets.includes(text)) {\n    // Set strikethrough on the paragraph's font (covers all runs) ...\n    paragraph.font.strikeThrough = true;\n    // ...and explicitly on the paragraph mark range so the trailing\n    // pilcrow carries the same rPr/strike as the other struck bullets.\n    const range = paragraph.getRange();\n    range.font.strikeThrough = true;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Apply strikethrough formatting to three specific \"deliverable\" bullet\n# paragraphs that were marked as done/resolved (data retention policy,\n# write/test triggers, and the \"at least two triggers must be for delete\"\n# bullet). The strike needs to land on every run AND on the paragraph\n# mark itself, matching how the other already-struck bullets in this\n# document are built.\n\n$targets = @(\n    \"Data retention policy: data is kept indefinitely. There are numerous ways to resolve this issue, find one and implement it.\",\n    \"Write and test seven (7) triggers for seven (7) separate tables to implement the business rules.\",\n    \"At least two (2) of the triggers must be for delete.\"\n)\n\n$d = $word.ActiveDocument\n\nforeach ($p in $d.Paragraphs) {\n    $text = $p.Range.Text.TrimEnd([char]13, [char]7)\n    if ($targets -contains $text) {\n        # Apply to the whole paragraph range (including the trailing\n        # paragraph mark) so both the runs and the pilcrow's rPr pick up\n        # the <w:strike/>.\n        $p.Range.Font.StrikeThrough = 1\n    }\n}\n"}
